{"js": "// Remove the trailing \"Ver no Jupiter...\" / footer copyright paragraphs (and\n// the blank paragraph right before them) that used to follow the\n// bibliography entry \"Rio de Janeiro: Elsevier Editora, 2007.\" while leaving\n// everything else (including the blank paragraph + page-break paragraph that\n// come after) untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"Rio de Janeiro: Elsevier Editora, 2007.\")\n// and the two text paragraphs that must be deleted, by matching their text\n// content rather than hard-coded indices, so the script is resilient to the\n// exact paragraph numbering.\nlet anchorIndex = -1;\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (text.indexOf(\"Rio de Janeiro: Elsevier Editora, 2007.\") !== -1) {\n    anchorIndex = i;\n  } else if (text.indexOf(\"Ver no Jupiter Salvar em pdf Salvar em docx\") !== -1) {\n    jupiterIndex = i;\n  } else if (text.indexOf(\"Contact: luizeleno@usp.br\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (anchorIndex === -1 || jupiterIndex === -1 || copyrightIndex === -1) {\n  throw new Error(\"Could not locate the expected paragraphs to delete.\");\n}\n\n// The blank paragraph that sits between the anchor paragraph and the\n// \"Ver no Jupiter...\" paragraph must go too.\nconst blankIndex = anchorIndex + 1;\nif (blankIndex >= jupiterIndex) {\n  throw new Error(\"Unexpected document layout around the anchor paragraph.\");\n}\n\n// Delete from the bottom up so earlier indices stay valid.\nitems[copyrightIndex].delete();\nitems[jupiterIndex].delete();\nitems[blankIndex].delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / footer copyright paragraphs (and\n# the blank paragraph right before them) that used to follow the\n# bibliography entry \"Rio de Janeiro: Elsevier Editora, 2007.\", leaving the\n# blank paragraph + page-break paragraph that come after untouched.\n\n$d = $word.ActiveDocument\n\n$anchorIndex = -1\n$jupiterIndex = -1\n$copyrightIndex = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.Contains(\"Rio de Janeiro: Elsevier Editora, 2007.\")) {\n        $anchorIndex = $i\n    } elseif ($t.Contains(\"Ver no Jupiter Salvar em pdf Salvar em docx\")) {\n        $jupiterIndex = $i\n    } elseif ($t.Contains(\"Contact: luizeleno@usp.br\")) {\n        $copyrightIndex = $i\n    }\n}\n\nif ($anchorIndex -eq -1 -or $jupiterIndex -eq -1 -or $copyrightIndex -eq -1) {\n    throw \"Could not locate the expected paragraphs to delete.\"\n}\n\n$blankIndex = $anchorIndex + 1\nif ($blankIndex -ge $jupiterIndex) {\n    throw \"Unexpected document layout around the anchor paragraph.\"\n}\n\n# Delete the blank paragraph, the \"Ver no Jupiter...\" paragraph and the\n# copyright paragraph as a single contiguous range so the document is only\n# mutated once. Range runs from the start of the blank paragraph through to\n# the start of the paragraph right after the copyright paragraph, which\n# removes all three paragraph marks along with their text.\n$startRange = $d.Paragraphs.Item($blankIndex).Range.Start\n$endPara = $d.Paragraphs.Item($copyrightIndex + 1)\n$endRange = $endPara.Range.Start\n\n$range = $d.Range($startRange, $endRange)\n$range.Delete()\n"}
